$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.193.65"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.834.72"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.01%  "
$orig = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.93"
$ws.Range("D5").Style = $orig
$ws.Range("E5").Value = "  +0.77%  "
$orig = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6654"
$ws.Range("D6").Style = $orig
$ws.Range("E6").Value = "  -2.30%  "
$orig = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = $orig
$ws.Range("E7").Value = "  +0.03%  "
$orig = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07411"
$ws.Range("D8").Style = $orig
$ws.Range("E8").Value = "  -0.65%  "
$orig = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2932"
$ws.Range("D9").Style = $orig
$ws.Range("E9").Value = "  -1.76%  "
$orig = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.93"
$ws.Range("D10").Style = $orig
$ws.Range("E10").Value = "  -0.88%  "
$orig = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07755"
$ws.Range("D11").Style = $orig
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "1.838.87"
$ws.Range("E12").Value = "  -0.30%  "
$orig = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.997"
$ws.Range("D13").Style = $orig
$ws.Range("E13").Value = "  -0.32%  "
$orig = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6679"
$ws.Range("D14").Style = $orig
$ws.Range("E14").Value = "  -1.32%  "
$orig = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.93"
$ws.Range("D15").Style = $orig
$ws.Range("E15").Value = "  -4.32%  "
$orig = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.113"
$ws.Range("D16").Style = $orig
$ws.Range("E16").Value = "  -0.71%  "
$orig = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008357"
$ws.Range("D17").Style = $orig
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "29.198.61"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "2.088.91"
$ws.Range("E19").Value = "  +0.06%  "
$orig = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.09"
$ws.Range("D20").Style = $orig
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.10%  "
$orig = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.152"
$ws.Range("D23").Style = $orig
$ws.Range("E23").Value = "  -2.67%  "
$orig = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("D24").Style = $orig
$ws.Range("E24").Value = "  +0.04%  "
$orig = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.20"
$ws.Range("D25").Style = $orig
$ws.Range("E25").Value = "  -1.09%  "
$orig = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1407"
$ws.Range("D26").Style = $orig
$ws.Range("E26").Value = "  -2.37%  "
$orig = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.615"
$ws.Range("D27").Style = $orig
$ws.Range("E27").Value = "  -1.07%  "
$orig = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.98"
$ws.Range("D28").Style = $orig
$ws.Range("E28").Value = "  -0.23%  "
$orig = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").Style = $orig
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  -3.28%  "
$orig = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.043"
$ws.Range("D31").Style = $orig
$ws.Range("E31").Value = "  -2.09%  "
$orig = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05288"
$ws.Range("D33").Style = $orig
$ws.Range("E33").Value = "  -2.33%  "
$orig = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.864"
$ws.Range("D34").Style = $orig
$ws.Range("E34").Value = "  +0.57%  "
$orig = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7456"
$ws.Range("D35").Style = $orig
$ws.Range("E35").Value = "  -0.98%  "
$orig = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.138"
$ws.Range("D36").Style = $orig
$ws.Range("E36").Value = "  +0.83%  "
$orig = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.647"
$ws.Range("D37").Style = $orig
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").Value = "1.293.04"
$ws.Range("E38").Value = "  -1.00%  "
$orig = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01799"
$ws.Range("D39").Style = $orig
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +0.71%  "
$orig = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9362"
$ws.Range("D41").Style = $orig
$ws.Range("E41").Value = "  +0.10%  "
$orig = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.881"
$ws.Range("D42").Style = $orig
$ws.Range("E42").Value = "  -2.95%  "
$orig = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08380"
$ws.Range("D43").Style = $orig
$ws.Range("E43").Value = "  -2.08%  "
$orig = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("D44").Style = $orig
$ws.Range("E44").Value = "  +0.10%  "
$orig = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.36"
$ws.Range("D45").Style = $orig
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").Value = "1.988.71"
$ws.Range("E46").Value = "  +0.28%  "
$orig = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5148"
$ws.Range("D47").Style = $orig
$ws.Range("E47").Value = "  -0.56%  "
$orig = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.759"
$ws.Range("D48").Style = $orig
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -0.68%  "
$orig = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "62.93"
$ws.Range("D50").Style = $orig
$ws.Range("E50").Value = "  -1.51%  "
$orig = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05879"
$ws.Range("D51").Style = $orig
$ws.Range("E51").Value = "  -0.77%  "
